$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed/last-updated) column C was refreshed from
# 2024-11-20 (45616) to 2024-11-21 (45617) for every data row (rows 2-33).
$ws.Range("C2:C33").Value = 45617
